$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.341.03"
$ws.Range("E2").Value = "  +4.50%  "
$ws.Range("D3").Value = "2.728.50"
$ws.Range("E3").Value = "  +3.63%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'529.71"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "'148.48"
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("D8").Value = "'0.581"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "2.752.09"
$ws.Range("E9").Value = "  +4.14%  "
$ws.Range("E10").Value = "  +13.24%  "
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  +2.73%  "
$ws.Range("E13").Value = "  +3.07%  "
$ws.Range("D14").Value = "3.202.80"
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").Value = "61.276.15"
$ws.Range("E15").Value = "  +4.37%  "
$ws.Range("E16").Value = "  +4.12%  "
$ws.Range("D17").Value = "2.737.10"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").Value = "'347.69"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("E21").Value = "  +3.89%  "
$ws.Range("D22").Value = "'6.45"
$ws.Range("E22").Value = "  +5.11%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'63.66"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("E25").Value = "  +5.17%  "
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "0.0₃0831"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("E29").Value = "  +5.34%  "
$ws.Range("D30").Value = "'6.76"
$ws.Range("E30").Value = "  +8.35%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").Value = "'19.15"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").Value = "'150.53"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("D35").Value = "'4.27"
$ws.Range("E35").Value = "  +7.54%  "
$ws.Range("E36").Value = "  +8.56%  "
$ws.Range("E37").Value = "  -6.57%  "
$ws.Range("E38").Value = "  +8.67%  "
$ws.Range("E39").Value = "  +9.22%  "
$ws.Range("D40").Value = "'37.49"
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("E43").Value = "  +4.29%  "
$ws.Range("D44").Value = "'282.35"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0989"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Value = "2.118.76"
$ws.Range("E48").Value = "  +7.70%  "
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").Value = "'19.57"
$ws.Range("E50").Value = "  +6.55%  "
$ws.Range("E51").Value = "  +2.24%  "

# Reset style on cells that required a quote-prefix to stay text,
# so no stray cell-level style reference is left behind.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").Style = "Normal"
